# POC- Retrieve Values from MDS File for Temperature
#
# The workbook originally has two tabs:
#   "Login"      - admin-style credentials (amrendrasadmin@yopmail.com / pass1234)
#   "LinksLogin" - link-based credentials (raviuser@yopmail.com / 12345678), active tab
#
# This edit repurposes the tabs:
#   - the sheet that holds the "raviuser" link-login data becomes the new,
#     first, active "Login" sheet
#   - the sheet that holds the old admin credentials is renamed "AdminLogin"
#     and pushed to the second position
#
# Worksheet object references in this host behave positionally once a Move()
# re-orders the tab strip, so sheets are re-fetched by name after each
# mutating step rather than reusing stale handles.

$wb = $excel.ActiveWorkbook

$oldLogin = $wb.Worksheets.Item("Login")
$oldLinks = $wb.Worksheets.Item("LinksLogin")

# Rename first (avoid a transient name collision) -- the old "Login" sheet
# becomes "AdminLogin", the old "LinksLogin" sheet becomes the new "Login".
$oldLogin.Name = "AdminLogin"
$oldLinks.Name = "Login"

# Re-fetch by the new names and move the new "Login" sheet in front of
# "AdminLogin" so the tab order is Login, AdminLogin.
$newLogin = $wb.Worksheets.Item("Login")
$newAdmin = $wb.Worksheets.Item("AdminLogin")
$newLogin.Move($newAdmin)

# Make "Login" the active/selected tab again and set each sheet's
# remembered selection.
$newLogin = $wb.Worksheets.Item("Login")
$newLogin.Activate()
$newLogin.Range("B18").Select()

$newAdmin = $wb.Worksheets.Item("AdminLogin")
$newAdmin.Range("A19").Select()

# Leave "Login" as the active sheet/selection when done.
$newLogin = $wb.Worksheets.Item("Login")
$newLogin.Activate()
